$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Column B width: 21.125 -> 22 (no longer "best fit", explicit custom width)
# ---------------------------------------------------------------------------
$ws.Columns("B").ColumnWidth = 21.125

# ---------------------------------------------------------------------------
# Build the "Blurbs" table skeleton (title + PK marker + header row) first.
# ---------------------------------------------------------------------------
$ws.Range("A44").Copy()
$ws.Range("A68").PasteSpecial(-4122)
$ws.Range("A68").Value = "Blurbs"

$ws.Range("A2").Copy()
$ws.Range("A69").PasteSpecial(-4122)
$ws.Range("A69").Value = "PK"

$ws.Range("A46:C46").Copy()
$ws.Range("A70:C70").PasteSpecial(-4122)
$ws.Range("A70").Value = "Student email"
$ws.Range("B70").Value = "Blurb"
$ws.Range("C70").Value = "Administrator email"

# ---------------------------------------------------------------------------
# Existing "Student Points" Status column: "yes" -> "accepted"
# ---------------------------------------------------------------------------
$ws.Range("C47").Value = "accepted"
$ws.Range("C48").Value = "accepted"
$ws.Range("C49").Value = "accepted"
$ws.Range("C53").Value = "accepted"
$ws.Range("C54").Value = "accepted"
$ws.Range("C57").Value = "accepted"

# ---------------------------------------------------------------------------
# Build the "Reflections" table skeleton (title + PK markers + header row).
# ---------------------------------------------------------------------------
$ws.Range("A44").Copy()
$ws.Range("A60").PasteSpecial(-4122)
$ws.Range("A60").Value = "Reflections"

$ws.Range("A45:B45").Copy()
$ws.Range("A61:B61").PasteSpecial(-4122)
$ws.Range("A61").Value = "PK"
$ws.Range("B61").Value = "PK"

$ws.Range("A46:B46").Copy()
$ws.Range("A62:B62").PasteSpecial(-4122)
$ws.Range("A62").Value = "StudentID"
$ws.Range("B62").Value = "Reflection"

# ---------------------------------------------------------------------------
# Fill in the "Blurbs" data rows (emails carry a hyperlink, like the other
# email columns elsewhere in the workbook).
# ---------------------------------------------------------------------------
$ws.Range("A71").Value = "rtravis@radford.edu"
$ws.Range("B71").Value = "Can you check my c2 please?"
$ws.Range("C71").Value = "pmartinez@radford.edu"

$ws.Range("A72").Value = "ljohnson@radford.edu"
$ws.Range("B72").Value = "When are the points for this semester due?"
$ws.Range("C72").Value = "othalwitz@radford.edu"

$ws.Range("A73").Value = "sjenkins@radford.edu"
$ws.Range("B73").Value = "How many points do I need?"
$ws.Range("C73").Value = "rhowards@radford.edu"

$ws.Range("A74").Value = "jmasterson@radford.edu"
$ws.Range("B74").Value = "Where can I find resources on how to get an internship?"
$ws.Range("C74").Value = "rhowards@radford.edu"

$ws.Hyperlinks.Add($ws.Range("A71"), "mailto:rtravis@radford.edu")
$ws.Hyperlinks.Add($ws.Range("A72"), "mailto:ljohnson@radford.edu")
$ws.Hyperlinks.Add($ws.Range("A73"), "mailto:sjenkins@radford.edu")
$ws.Hyperlinks.Add($ws.Range("A74"), "mailto:jmasterson@radford.edu")
$ws.Hyperlinks.Add($ws.Range("C71"), "mailto:pmartinez@radford.edu")
$ws.Hyperlinks.Add($ws.Range("C72"), "mailto:othalwitz@radford.edu")
$ws.Hyperlinks.Add($ws.Range("C73"), "mailto:rhowards@radford.edu")
$ws.Hyperlinks.Add($ws.Range("C74"), "mailto:rhowards@radford.edu")

# Restore the plain "hyperlink-styled" formatting (Hyperlinks.Add pushes its
# own style onto the cell) so the cells match the look of D4 (the existing
# hyperlink cells elsewhere in the sheet).
$ws.Range("D4").Copy()
$ws.Range("A71:A74").PasteSpecial(-4122)
$ws.Range("C71:C74").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Fill in the "Reflections" data rows.
# ---------------------------------------------------------------------------
$ws.Range("A63").Value = 900752513
$ws.Range("B63").Value = "refelction text"

$ws.Range("A64").Value = 922152345
$ws.Range("B64").Value = "Reflection text"

$ws.Range("A65").Value = 900019812
$ws.Range("B65").Value = "reflection text"

$ws.Range("A66").Value = 956789000
$ws.Range("B66").Value = "reflection text"

# ---------------------------------------------------------------------------
# Final view state: scroll/select like the author left the sheet.
# ---------------------------------------------------------------------------
$ws.Range("C67").Select()
